$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.504.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.21"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.63%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.88"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07544"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.2970"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.23"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07680"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.886.57"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6846"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000009800"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.128.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.210"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "29.546.68"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "233.96"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.609"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.76"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1388"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.432"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05840"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.282"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.108"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.891"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.172"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.7159"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.591"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.800"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.235.83"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9124"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.127"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.040.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.93"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.46"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.286"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.162"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.75%  "
